$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-04 Thursday", "2025-09-05 Friday"),
    @("518÷5=", "610÷6="),
    @("871÷3=", "232÷8="),
    @("816÷4=", "743÷6="),
    @("750÷7=", "989÷7="),
    @("648÷4=", "895÷7="),
    @("356÷8=", "625÷5="),
    @("459÷6=", "663÷2="),
    @("704÷7=", "640÷4="),
    @("154÷5=", "135÷8="),
    @("755÷8=", "910÷6="),
    @("887÷5=", "140÷9="),
    @("506÷4=", "242÷6="),
    @("491÷9=", "653÷6="),
    @("721÷8=", "951÷8="),
    @("638÷5=", "826÷4="),
    @("668÷4=", "318÷6="),
    @("703÷3=", "230÷3="),
    @("505÷7=", "102÷9="),
    @("466÷5=", "695÷9="),
    @("262÷3=", "793÷8="),
    @("637÷6=", "863÷4="),
    @("562÷2=", "732÷4="),
    @("482÷3=", "505÷8="),
    @("815÷3=", "567÷5="),
    @("164÷9=", "495÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
